# correccion en la generacion de comprobante
# Actualiza el recibo de pago con el destinatario, concepto y numero de
# comprobante correctos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numero de comprobante (P3): N°:248 -> N°:267
$ws.Range("P3").Value = "N°:267"

# Destinatario + DNI (L7): SR N , DNI 111 -> ALVAREZ MARIA ALISA , DNI 33051226
$ws.Range("L7").Value = "ALVAREZ MARIA ALISA , DNI 33051226"

# Concepto (C8): Cuota3 -> Cuota1
$ws.Range("C8").Value = "Cuota1"

# Domicilio (L9): A1 -> 20 (debe seguir siendo texto, no numero, y
# conservar el estilo original de la celda; se usa una celda auxiliar
# fuera del area usada de la hoja + copiar-y-pegar-valores para lograrlo
# sin alterar el formato numerico de L9).
$scratch = $ws.Range("ZZ500")
$scratch.Value = "'20"
$scratch.Copy()
$ws.Range("L9").PasteSpecial(-4163)
$scratch.Clear()

# Detalle (I13): En concepto de pago en efectivo por Cuota3 -> ... Cuota1
$ws.Range("I13").Value = "En concepto de pago en efectivo por Cuota1"

# Observaciones (K18) permanece vacio
$ws.Range("K18").Value = ""
